$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, matching the header style used by the
# other header cells (bold, centered, bordered) via a format-only copy
# from the neighboring header cell G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data value 0 for the "Save" column on row 2
$ws.Range("H2").Value = 0
